# docs/validacao.xlsx - "testes debug e agent"
# Adds two new SBERT model comparison columns (G, H) and repurposes the
# existing GT / SBERT columns, filling in pass/fail results for every row
# of the validation matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$okOk     = "OK na busca - OK no LLM"
$failFail = "Fail na busca - Fail no LLM"

# ---- Column F: SBERT distiluse ------------------------------------------
$ws.Cells.Item(2, 6).Value = "SBERT distiluse"
$ws.Cells.Item(3, 6).Value = $failFail
$ws.Cells.Item(4, 6).Value = $okOk
$ws.Cells.Item(5, 6).Value = $failFail
$ws.Cells.Item(6, 6).Value = $okOk
$ws.Cells.Item(7, 6).Value = $failFail
$ws.Cells.Item(8, 6).Value = $failFail

# ---- Column G: SBERT paraphrase MiniLM -----------------------------------
$ws.Cells.Item(2, 7).Value = "SBERT paraphrase MiniLM"
$ws.Cells.Item(3, 7).Value = $failFail
$ws.Cells.Item(4, 7).Value = $okOk
$ws.Cells.Item(5, 7).Value = $failFail
$ws.Cells.Item(6, 7).Value = $okOk
$ws.Cells.Item(7, 7).Value = $failFail
$ws.Cells.Item(8, 7).Value = $okOk

# ---- Column H: SBERT paraphrase mpnet ------------------------------------
$ws.Cells.Item(2, 8).Value = "SBERT paraphrase mpnet"
$ws.Cells.Item(3, 8).Value = $failFail
$ws.Cells.Item(4, 8).Value = $okOk
$ws.Cells.Item(5, 8).Value = $failFail
$ws.Cells.Item(6, 8).Value = $failFail
$ws.Cells.Item(7, 8).Value = $failFail
$ws.Cells.Item(8, 8).Value = $okOk

# ---- Column D header: GT -> Ground Truth (renamed last) ------------------
$ws.Cells.Item(2, 4).Value = "Ground Truth"

# ---- Column widths for the newly-used / widened columns ------------------
$ws.Columns.Item(4).ColumnWidth = 15.67   # D: widen for "Ground Truth"
$ws.Columns.Item(6).ColumnWidth = 25.5025 # F
$ws.Columns.Item(7).ColumnWidth = 25.17   # G
$ws.Columns.Item(8).ColumnWidth = 28.5025 # H

# ---- Selection matches the saved cursor position in the author's file ----
$ws.Range("G10").Select()
